$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -10.96189999999999
$ws.Range("C10").Value = -12.4152
$ws.Range("C12").Value = -14.4007
$ws.Range("C18").Value = -14.248
$ws.Range("C25").Value = -10.85839999999999
